# Update for next week.
$wb = $excel.ActiveWorkbook
$zadania = $wb.Worksheets.Item("Zadania")
$week = $wb.Worksheets.Item("19-25.10")

# --- Sheet "19-25.10": new task "Bazy danych" (Praca magisterska) added to the sidebar list ---
$week.Range("A13").Value = "Bazy danych"

# --- Sheet "Zadania": add the MGR (Praca magisterska) row ---
$zadania.Range("A12").Value = "Praca magisterska"
$zadania.Range("F12").Value = "MGR"

# --- Sheet "Zadania": roll the September dates forward to October ---
$zadania.Range("D6").Value = "23.10.2014"
$zadania.Range("D7").Value = "20.10.2014"
$zadania.Range("D8").Value = "21.10.2014"
$zadania.Range("D9").Value = "22.10.2014"
$zadania.Range("D10").Value = "23.10.2014"

# --- Sheet "19-25.10": weekly planner updates ---
$week.Range("D21").Value = "dokończ DI 1h"
$week.Range("E21").Value = "PF - 1,5h"
$week.Range("H22").Value = "SU "
$week.Range("C23").Value = "Mail MGR"
$week.Columns.Item(4).ColumnWidth = 15.71

# --- Restore the author's last-saved selections ---
$zadania.Range("E14").Select()
$week.Range("C24").Select()
